$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.096.06"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +3.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.988.97"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.84%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.32"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.04"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +10.37%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.980.53"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.133"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +6.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.11"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +8.14%  "
$ws.Range("E12").Value = "  +4.47%  "
$ws.Range("E13").Value = "  +7.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.67"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.78%  "
$ws.Range("E15").Value = "  +2.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.488.47"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.08"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +7.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.989.69"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "59.077.26"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "428.11"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +5.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.57"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +5.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.714"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +6.48%  "
$ws.Range("E23").Value = "  +3.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.42"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +5.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.54"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.14%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("E28").Value = "  +10.93%  "
$ws.Range("E29").Value = "  +3.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.74"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +6.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.71"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.12"
$ws.Range("D32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0985"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("E34").Value = "  +21.97%  "
$ws.Range("E35").Value = "  +8.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.88"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +7.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.07"
$ws.Range("D37").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.32"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.62"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.72"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +11.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "401.40"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +9.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.780.01"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0351"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.73%  "
$ws.Range("E44").Value = "  +0.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.252"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +10.11%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.49"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.42%  "
$ws.Range("E48").Value = "  +2.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.72"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +21.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.00"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.42"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.60%  "
